$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16/17: swap WrappedEther <-> ShibaInu (name, link), with new price/volume ---
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"

# --- Column D (Price) updates that are safe as plain text (not numeric-parseable) ---
$ws.Range("D2").Value = "58.782.12"
$ws.Range("D3").Value = "2.650.78"
$ws.Range("D13").Value = "3.116.67"
$ws.Range("D14").Value = "58.830.92"
$ws.Range("D17").Value = "2.641.08"
$ws.Range("D27").Value = "0.0₃0798"
$ws.Range("D48").Value = "2.035.09"

# --- Column D (Price) updates that look numeric: force text format first so Excel keeps literal text (e.g. "1.00") ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D5").Value = "523.07"
$ws.Range("D6").Value = "144.11"
$ws.Range("D7").Value = "1.00"
$ws.Range("D9").Value = "7.00"
$ws.Range("D10").Value = "0.102"
$ws.Range("D11").Value = "0.333"
$ws.Range("D15").Value = "21.02"
$ws.Range("D16").Value = "0.0000136"
$ws.Range("D18").Value = "338.67"
$ws.Range("D19").Value = "4.35"
$ws.Range("D20").Value = "10.34"
$ws.Range("D21").Value = "6.34"
$ws.Range("D22").Value = "1.00"
$ws.Range("D23").Value = "63.77"
$ws.Range("D24").Value = "0.417"
$ws.Range("D25").Value = "0.165"
$ws.Range("D26").Value = "1.00"
$ws.Range("D28").Value = "7.07"
$ws.Range("D29").Value = "6.68"
$ws.Range("D30").Value = "0.999"
$ws.Range("D32").Value = "18.81"
$ws.Range("D33").Value = "149.59"
$ws.Range("D34").Value = "4.14"
$ws.Range("D36").Value = "0.887"
$ws.Range("D37").Value = "0.867"
$ws.Range("D38").Value = "36.71"
$ws.Range("D40").Value = "3.58"
$ws.Range("D41").Value = "0.615"
$ws.Range("D42").Value = "0.999"
$ws.Range("D43").Value = "19.88"
$ws.Range("D44").Value = "275.34"
$ws.Range("D45").Value = "0.0967"
$ws.Range("D47").Value = "0.0532"
$ws.Range("D49").Value = "4.71"
$ws.Range("D50").Value = "0.0227"
$ws.Range("D51").Value = "18.83"

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = "  -2.52%  "
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("E6").Value = "  -1.05%  "
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("E8").Value = "  -1.54%  "
$ws.Range("E9").Value = "  +8.74%  "
$ws.Range("E10").Value = "  -2.52%  "
$ws.Range("E11").Value = "  -1.96%  "
$ws.Range("E12").Value = "  +1.48%  "
$ws.Range("E13").Value = "  -0.87%  "
$ws.Range("E14").Value = "  -2.46%  "
$ws.Range("E15").Value = "  -1.17%  "
$ws.Range("E16").Value = "  -1.81%  "
$ws.Range("E17").Value = "  -4.40%  "
$ws.Range("E18").Value = "  -3.22%  "
$ws.Range("E19").Value = "  -4.27%  "
$ws.Range("E20").Value = "  -1.78%  "
$ws.Range("E21").Value = "  +0.60%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("E23").Value = "  +1.72%  "
$ws.Range("E24").Value = "  -0.86%  "
$ws.Range("E25").Value = "  -1.48%  "
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("E27").Value = "  -1.70%  "
$ws.Range("E28").Value = "  -2.40%  "
$ws.Range("E29").Value = "  -2.63%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("E32").Value = "  -1.30%  "
$ws.Range("E33").Value = "  +0.85%  "
$ws.Range("E34").Value = "  -3.63%  "
$ws.Range("E35").Value = "  -2.88%  "
$ws.Range("E36").Value = "  -6.56%  "
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("E39").Value = "  -6.10%  "
$ws.Range("E40").Value = "  -3.06%  "
$ws.Range("E41").Value = "  +1.06%  "
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("E44").Value = "  -1.89%  "
$ws.Range("E45").Value = "  -2.28%  "
$ws.Range("E46").Value = "  +1.97%  "
$ws.Range("E47").Value = "  -1.52%  "
$ws.Range("E48").Value = "  -4.24%  "
$ws.Range("E49").Value = "  -2.72%  "
$ws.Range("E50").Value = "  -2.82%  "
$ws.Range("E51").Value = "  -1.00%  "
